$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 46 (2000, Group B, Turkey)
$ws.Range("F46").Value = 0
$ws.Range("H46").Value = 1

# Row 48 (2000, Group B, Sweden)
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0

# Rows 149-152 (2021, Group F) tie-break reshuffle
$ws.Range("C149").Value = "Portugal"
$ws.Range("D149").Value = 2
$ws.Range("H149").Value = 3

$ws.Range("C150").Value = "Germany"
$ws.Range("D150").Value = 0
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 2

$ws.Range("C151").Value = "France"
$ws.Range("D151").Value = 1
$ws.Range("E151").Value = 3
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 6

$ws.Range("E152").Value = 0
$ws.Range("F152").Value = 2
$ws.Range("H152").Value = 4

# Rows 162-164 (2024, Group C)
$ws.Range("E162").Value = 1
$ws.Range("H162").Value = 2

$ws.Range("F163").Value = 1
$ws.Range("H163").Value = 2

$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 0
